# "Doing Updates for Financials"
# Update the ASPU yearly-financials figures on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non Recurring (row 14) - clear out the old mixed NA/values, set all to 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0

# Total Operating Expenses (row 17)
$ws.Range("D17").Value = 27300
$ws.Range("G17").Value = 8600

# Operating Income or Loss (row 18)
$ws.Range("D18").Value = -5300
$ws.Range("G18").Value = -3400

# Total Other Income/Expenses Net (row 20)
$ws.Range("D20").Value = 200
$ws.Range("G20").Value = -400

# Other Items (row 32)
$ws.Range("D32").Value = -200
$ws.Range("G32").Value = 400
